$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the betting-data columns (B, E:AD) between row 18 and row 19
$rowA = $ws.Range("B18:AD18")
$rowB = $ws.Range("B19:AD19")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

# Swap the betting-data columns (B, E:AD) between row 45 and row 46
$rowA = $ws.Range("B45:AD45")
$rowB = $ws.Range("B46:AD46")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

# Swap the betting-data columns (B, E:AD) between row 50 and row 51
$rowA = $ws.Range("B50:AD50")
$rowB = $ws.Range("B51:AD51")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

# Swap the betting-data columns (B, E:AD) between row 62 and row 63
$rowA = $ws.Range("B62:AD62")
$rowB = $ws.Range("B63:AD63")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

# Swap the betting-data columns (B, E:AD) between row 96 and row 97
$rowA = $ws.Range("B96:AD96")
$rowB = $ws.Range("B97:AD97")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

# Swap the betting-data columns (B, E:AD) between row 102 and row 103
$rowA = $ws.Range("B102:AD102")
$rowB = $ws.Range("B103:AD103")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

# Swap the betting-data columns (B, E:AD) between row 118 and row 119
$rowA = $ws.Range("B118:AD118")
$rowB = $ws.Range("B119:AD119")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

# Swap the betting-data columns (B, E:AD) between row 129 and row 130
$rowA = $ws.Range("B129:AD129")
$rowB = $ws.Range("B130:AD130")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

# Swap the betting-data columns (B, E:AD) between row 131 and row 132
$rowA = $ws.Range("B131:AD131")
$rowB = $ws.Range("B132:AD132")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA
